# "for TOT training session"
# Update the October data-collection window: startDate (D2) moves from
# the 2nd to the 18th, and endDate (E2) moves from the 6th to the 19th.
# no_of_days (F2) recalculates automatically via its existing formula
# (=E2-D2+1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 18
$ws.Range("E2").Value = 19

# Move the active selection to E3, matching where the cursor was left
# after entering the new endDate value.
$ws.Range("E3").Select()
